$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Day 3): add squat, deadlift, row goal values
$ws.Range("C6").Value = 350
$ws.Range("D6").Value = 450
$ws.Range("E6").Value = 150

# Row 7 (Day 4): add benchpress, squat, deadlift, row goal values
$ws.Range("B7").Value = 210
$ws.Range("C7").Value = 360
$ws.Range("D7").Value = 450
$ws.Range("E7").Value = 150
